$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Workbook-level: rename Sheet2 -> "Emails"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)   # Users
$ws2 = $wb.Worksheets.Item(2)   # Sheet2 -> Emails
$ws3 = $wb.Worksheets.Item(3)   # Sheet3

$ws2.Name = "Emails"

# ---------------------------------------------------------------------------
# 2. Users sheet: append the new Cobalt user rows (53-80)
# ---------------------------------------------------------------------------
$newUsers = @(
    @{Row=53; User='SearchOpenWebUser1'; Email='SearchOpenWeb@mailinator.com '}
    @{Row=54; User='FFHUser1'; Email='FFHUser1@mailinator.com '}
    @{Row=55; User='FFHUser2'; Email='FFHUser2@mailinator.com'}
    @{Row=56; User='FFHUser3'; Email='FFHUser3@mailinator.com'}
    @{Row=57; User='FFHUser4'; Email='FFHUser4@mailinator.com'}
    @{Row=58; User='FrontEndUser1'; Email='FrontEndUser1@mailinator.com'}
    @{Row=59; User='FrontEndUser2'; Email='FrontEndUser2@mailinator.com'}
    @{Row=60; User='FrontEndUser3'; Email='FrontEndUser3@mailinator.com'}
    @{Row=61; User='FrontEndUser4'; Email='FrontEndUser4@mailinator.com'}
    @{Row=62; User='FrontEndUser5'; Email='FrontEndUser5@mailinator.com'}
    @{Row=63; User='FrontEndUser6'; Email='FrontEndUser6@mailinator.com'}
    @{Row=64; User='FrontEndUser7'; Email='FrontEndUser7@mailinator.com'}
    @{Row=65; User='FrontEndUser8'; Email='FrontEndUser8@mailinator.com'}
    @{Row=66; User='FrontEndUser9'; Email='FrontEndUser9@mailinator.com'}
    @{Row=67; User='FrontEndUser10'; Email='FrontEndUser10@mailinator.com'}
    @{Row=68; User='UrlUser1'; Email='UrlUser1@mailinator.com'}
    @{Row=69; User='UrlUser2'; Email='UrlUser2@mailinator.com'}
    @{Row=70; User='UrlUser3'; Email='UrlUser3@mailinator.com'}
    @{Row=71; User='LinkingUser1'; Email='LinkingUser1@mailinator.com '}
    @{Row=72; User='LoginUser1'; Email='LoginUser1@mailinator.com '}
    @{Row=73; User='LoginUser2'; Email='LoginUser2@mailinator.com '}
    @{Row=74; User='LoginUser3'; Email='LoginUser3@mailinator.com '}
    @{Row=75; User='LoginUser4'; Email='LoginUser4@mailinator.com '}
    @{Row=76; User='LoginUser5'; Email='LoginUser5@mailinator.com '}
    @{Row=77; User='LoginUser6'; Email='LoginUser6@mailinator.com '}
    @{Row=78; User='LoginUser7'; Email='LoginUser7@mailinator.com '}
    @{Row=79; User='CpetUser1'; Email='CpetUser1@mailinator.com '}
    @{Row=80; User='CpetUser2'; Email='CpetUser2@mailinator.com '}
)

foreach ($item in $newUsers) {
    $r = $item.Row

    $ws1.Cells.Item($r, 1).Value = $item.User
    $ws1.Cells.Item($r, 2).Value = "Password1"

    $descCell = $ws1.Cells.Item($r, 5)
    $descCell.Value = "THIS IS IN USE 24/7 - DO NOT USE!"
    $descCell.Borders(10).LineStyle = 1
    $descCell.Borders(10).Weight = 2
    $descCell.Borders(7).LineStyle = 1
    $descCell.Borders(7).Weight = 2

    $lockCell = $ws1.Cells.Item($r, 6)
    $lockCell.Value = "N"
    $lockCell.Borders(10).LineStyle = 1
    $lockCell.Borders(10).Weight = 2
    $lockCell.Borders(7).LineStyle = 1
    $lockCell.Borders(7).Weight = 2

    $emailCell = $ws1.Cells.Item($r, 7)
    $emailCell.Value = $item.Email
    $ws1.Hyperlinks.Add($emailCell, "mailto:" + $item.Email)
}

# Four trailing spacer rows (81-84) - just the bordered "E" cell, no value
for ($r = 81; $r -le 84; $r++) {
    $spacerCell = $ws1.Cells.Item($r, 5)
    $spacerCell.Borders(7).LineStyle = 1
    $spacerCell.Borders(7).Weight = 2
    $spacerCell.Borders(10).LineStyle = 1
    $spacerCell.Borders(10).Weight = 2
}

# ---------------------------------------------------------------------------
# 3. Column widths / view tidy-up on the Users sheet
# ---------------------------------------------------------------------------
$ws1.Columns.Item(1).AutoFit()
$ws1.Columns.Item(5).AutoFit()

# ---------------------------------------------------------------------------
# 4. Emails sheet: new lookup table
# ---------------------------------------------------------------------------
$ws2.Cells.Item(1, 1).Value = "Email"
$ws2.Cells.Item(1, 2).Value = "Password"
$ws2.Cells.Item(2, 1).Value = "tr-anz-tester1@yandex.com"
$ws2.Cells.Item(2, 2).Value = "tranztest"
$ws2.Cells.Item(3, 1).Value = "tr-anz-tester2@yandex.com"
$ws2.Cells.Item(3, 2).Value = "tranztest"

$ws2.Columns.Item(1).AutoFit()
$ws2.Range("A1:B3").Select()

# ---------------------------------------------------------------------------
# 5. Sheet3: drop the duplicate rows 3-5, keep only row 2
# ---------------------------------------------------------------------------
$ws3.Rows.Item(3).Delete()
$ws3.Rows.Item(3).Delete()
$ws3.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# 6. Restore Users as the active sheet / selection, matching the authored file
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("C80").Select()

